# Auto-generated Excel COM-interop script to apply crypto price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value while forcing text storage (avoids Excel auto-converting
# numeric-looking strings like "1.00" or "0.999" into real numbers, which would drop
# formatting such as trailing zeros). Using a temporary Text number format plus a
# ClearFormats() afterwards keeps the cell style untouched (no "s" attribute added).
function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

$ws.Range("D2").Value = '79.175.08'
$ws.Range("E2").Value = '  +3.75%  '
$ws.Range("D3").Value = '3.140.78'
$ws.Range("E3").Value = '  +1.75%  '
$ws.Range("E4").Value = '  +0.07%  '
Set-TextValue $ws.Range("D5") '204.03'
$ws.Range("E5").Value = '  +2.78%  '
Set-TextValue $ws.Range("D6") '621.59'
$ws.Range("E6").Value = '  +0.47%  '
Set-TextValue $ws.Range("D7") '0.260'
$ws.Range("E7").Value = '  +24.41%  '
Set-TextValue $ws.Range("D8") '1.00'
$ws.Range("E8").Value = '  +0.03%  '
Set-TextValue $ws.Range("D9") '0.584'
$ws.Range("E9").Value = '  +5.89%  '
$ws.Range("D10").Value = '3.141.10'
$ws.Range("E10").Value = '  +1.78%  '
Set-TextValue $ws.Range("D11") '0.580'
$ws.Range("E11").Value = '  +31.65%  '
Set-TextValue $ws.Range("D12") '0.0000248'
$ws.Range("E12").Value = '  +27.51%  '
$ws.Range("E13").Value = '  +1.76%  '
$ws.Range("D14").Value = '3.725.44'
$ws.Range("E14").Value = '  +2.01%  '
Set-TextValue $ws.Range("D15") '5.23'
$ws.Range("E15").Value = '  -0.24%  '
Set-TextValue $ws.Range("D16") '31.02'
$ws.Range("E16").Value = '  +5.88%  '
$ws.Range("D17").Value = '79.378.05'
$ws.Range("E17").Value = '  +4.17%  '
$ws.Range("D18").Value = '3.148.61'
$ws.Range("E18").Value = '  +2.26%  '
Set-TextValue $ws.Range("D19") '14.14'
$ws.Range("E19").Value = '  +4.13%  '
Set-TextValue $ws.Range("D20") '432.90'
$ws.Range("E20").Value = '  +13.31%  '
Set-TextValue $ws.Range("D21") '8.99'
$ws.Range("E21").Value = '  -0.31%  '
Set-TextValue $ws.Range("D22") '2.88'
$ws.Range("E22").Value = '  +10.78%  '
Set-TextValue $ws.Range("D23") '5.15'
$ws.Range("E23").Value = '  +15.82%  '
Set-TextValue $ws.Range("D24") '6.75'
$ws.Range("E24").Value = '  +4.65%  '
$ws.Range("D25").Value = '3.325.79'
$ws.Range("E25").Value = '  +2.65%  '
Set-TextValue $ws.Range("D26") '75.53'
$ws.Range("E26").Value = '  +4.32%  '
Set-TextValue $ws.Range("D27") '4.61'
$ws.Range("E27").Value = '  +2.90%  '
Set-TextValue $ws.Range("D28") '10.66'
$ws.Range("E28").Value = '  +5.73%  '
$ws.Range("E29").Value = '  -0.29%  '
Set-TextValue $ws.Range("D30") '0.0000119'
$ws.Range("E30").Value = '  +9.54%  '
Set-TextValue $ws.Range("D31") '0.999'
$ws.Range("E31").Value = '  +0.12%  '
Set-TextValue $ws.Range("D32") '8.85'
$ws.Range("E32").Value = '  +6.42%  '
Set-TextValue $ws.Range("D33") '541.68'
$ws.Range("E33").Value = '  +7.73%  '
Set-TextValue $ws.Range("D34") '1.46'
$ws.Range("E34").Value = '  +1.99%  '
Set-TextValue $ws.Range("D35") '1.98'
$ws.Range("E35").Value = '  +2.69%  '
Set-TextValue $ws.Range("D36") '0.147'
$ws.Range("E36").Value = '  +18.77%  '
Set-TextValue $ws.Range("D37") '22.73'
$ws.Range("E37").Value = '  +9.19%  '
Set-TextValue $ws.Range("D38") '0.122'
$ws.Range("E38").Value = '  +18.96%  '
$ws.Range("E39").Value = '  -0.07%  '
Set-TextValue $ws.Range("D40") '0.401'
$ws.Range("E40").Value = '  +5.69%  '
Set-TextValue $ws.Range("D41") '20.71'
$ws.Range("E41").Value = '  +3.22%  '
Set-TextValue $ws.Range("D42") '162.59'
$ws.Range("E42").Value = '  +0.47%  '
$ws.Range("B43").Value = 'USDe'
$ws.Range("C43").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextValue $ws.Range("D43") '1.00'
$ws.Range("E43").Value = '  +0.01%  '
$ws.Range("B44").Value = 'RenderToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
Set-TextValue $ws.Range("D44") '5.52'
$ws.Range("E44").Value = '  +7.25%  '
Set-TextValue $ws.Range("D45") '185.93'
$ws.Range("E45").Value = '  -5.05%  '
Set-TextValue $ws.Range("D46") '1.78'
$ws.Range("E46").Value = '  +7.42%  '
Set-TextValue $ws.Range("D47") '2.63'
$ws.Range("E47").Value = '  +8.15%  '
Set-TextValue $ws.Range("D48") '0.774'
$ws.Range("E48").Value = '  -4.05%  '
$ws.Range("B49").Value = 'OKB'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue $ws.Range("D49") '42.94'
$ws.Range("E49").Value = '  +4.09%  '
$ws.Range("B50").Value = 'ImmutableX'
$ws.Range("C50").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue $ws.Range("D50") '1.28'
$ws.Range("E50").Value = '  +2.44%  '
Set-TextValue $ws.Range("D51") '4.19'
$ws.Range("E51").Value = '  +7.16%  '
